# Adds the "ODI Bowling Extra" worksheet (sheetId 5) with its MATCH_CODE /
# MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL columns, mirroring the layout that
# already exists on the "ODI Batting Extra" sheet.

$wb = $excel.ActiveWorkbook

# --- create the new sheet, placed after the last existing tab -------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Bowling Extra"

# --- header row: copy the header formatting already used on the other
#     "Extra" sheet so the new header matches it exactly (bold, centered,
#     bordered), then fill in this sheet's own header text -----------------
$headerSource = $wb.Worksheets.Item("ODI Batting Extra").Range("A1")
$headerSource.Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)

$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "MAIDEN_OVERS"
$ws.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# --- data rows --------------------------------------------------------------
# MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL ($null = left blank)
$rows = @(
    @("4414", $null, $null),
    @("4417", "0",   "20.00%"),
    @("4449", "1",   "20.00%"),
    @("4450", $null, $null),
    @("4451", "0",   $null),
    @("4483", "0",   $null),
    @("4484", "0",   "10.00%"),
    @("4486", $null, $null),
    @("4519", $null, $null),
    @("4520", "1",   $null),
    @("4522", $null, $null),
    @("4533", $null, $null),
    @("4535", "2",   "10.00%"),
    @("4536", "1",   "40.00%"),
    @("4624", "0",   $null),
    @("4636", "1",   "20.00%"),
    @("4639", "2",   "30.00%"),
    @("4642", $null, $null),
    @("4727", $null, $null),
    @("4731", "1",   $null)
)

$r = 2
foreach ($row in $rows) {
    $matchCode   = $row[0]
    $maidenOvers = $row[1]
    $percentWkts = $row[2]

    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $matchCode

    if ($null -ne $maidenOvers) {
        $ws.Range("B$r").NumberFormat = "@"
        $ws.Range("B$r").Value = $maidenOvers
    }

    if ($null -ne $percentWkts) {
        $ws.Range("C$r").NumberFormat = "@"
        $ws.Range("C$r").Value = $percentWkts
    }

    $r++
}
